{"js": "const replacements = [\n  [\"2025-04-07 Monday\", \"2025-04-08 Tuesday\"],\n  [\"890\u00d76=5340\", \"178\u00d74=712\"],\n  [\"779\u00d74=3116\", \"726\u00d74=2904\"],\n  [\"555\u00d74=2220\", \"792\u00d72=1584\"],\n  [\"481\u00d75=2405\", \"155\u00d73=465\"],\n  [\"209\u00d75=1045\", \"607\u00d79=5463\"],\n  [\"324\u00d79=2916\", \"331\u00d79=2979\"],\n  [\"731\u00d75=3655\", \"719\u00d75=3595\"],\n  [\"377\u00d77=2639\", \"901\u00d74=3604\"],\n  [\"845\u00d74=3380\", \"901\u00d74=3604\"],\n  [\"297\u00d77=2079\", \"217\u00d72=434\"],\n  [\"595\u00d74=2380\", \"209\u00d76=1254\"],\n  [\"224\u00d72=448\", \"170\u00d78=1360\"],\n  [\"188\u00d75=940\", \"804\u00d75=4020\"],\n  [\"991\u00d74=3964\", \"837\u00d77=5859\"],\n  [\"450\u00d79=4050\", \"775\u00d73=2325\"],\n  [\"200\u00d77=1400\", \"557\u00d77=3899\"],\n  [\"745\u00d72=1490\", \"155\u00d76=930\"],\n  [\"757\u00d72=1514\", \"369\u00d78=2952\"],\n  [\"714\u00d74=2856\", \"627\u00d76=3762\"],\n  [\"249\u00d78=1992\", \"263\u00d74=1052\"],\n  [\"290\u00d75=1450\", \"669\u00d73=2007\"],\n  [\"739\u00d78=5912\", \"147\u00d74=588\"],\n  [\"759\u00d76=4554\", \"287\u00d78=2296\"],\n  [\"766\u00d76=4596\", \"723\u00d75=3615\"],\n  [\"374\u00d76=2244\", \"493\u00d77=3451\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "function Replace-Text($Document, $OldText, $NewText) {\n    $find = $Document.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $OldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $NewText\n    $find.Execute([ref]$OldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$NewText, 2) | Out-Null\n}\n\n$d = $word.ActiveDocument\n\nReplace-Text $d \"2025-04-07 Monday\" \"2025-04-08 Tuesday\"\nReplace-Text $d \"890\u00d76=5340\" \"178\u00d74=712\"\nReplace-Text $d \"779\u00d74=3116\" \"726\u00d74=2904\"\nReplace-Text $d \"555\u00d74=2220\" \"792\u00d72=1584\"\nReplace-Text $d \"481\u00d75=2405\" \"155\u00d73=465\"\nReplace-Text $d \"209\u00d75=1045\" \"607\u00d79=5463\"\nReplace-Text $d \"324\u00d79=2916\" \"331\u00d79=2979\"\nReplace-Text $d \"731\u00d75=3655\" \"719\u00d75=3595\"\nReplace-Text $d \"377\u00d77=2639\" \"901\u00d74=3604\"\nReplace-Text $d \"845\u00d74=3380\" \"901\u00d74=3604\"\nReplace-Text $d \"297\u00d77=2079\" \"217\u00d72=434\"\nReplace-Text $d \"595\u00d74=2380\" \"209\u00d76=1254\"\nReplace-Text $d \"224\u00d72=448\" \"170\u00d78=1360\"\nReplace-Text $d \"188\u00d75=940\" \"804\u00d75=4020\"\nReplace-Text $d \"991\u00d74=3964\" \"837\u00d77=5859\"\nReplace-Text $d \"450\u00d79=4050\" \"775\u00d73=2325\"\nReplace-Text $d \"200\u00d77=1400\" \"557\u00d77=3899\"\nReplace-Text $d \"745\u00d72=1490\" \"155\u00d76=930\"\nReplace-Text $d \"757\u00d72=1514\" \"369\u00d78=2952\"\nReplace-Text $d \"714\u00d74=2856\" \"627\u00d76=3762\"\nReplace-Text $d \"249\u00d78=1992\" \"263\u00d74=1052\"\nReplace-Text $d \"290\u00d75=1450\" \"669\u00d73=2007\"\nReplace-Text $d \"739\u00d78=5912\" \"147\u00d74=588\"\nReplace-Text $d \"759\u00d76=4554\" \"287\u00d78=2296\"\nReplace-Text $d \"766\u00d76=4596\" \"723\u00d75=3615\"\nReplace-Text $d \"374\u00d76=2244\" \"493\u00d77=3451\"\n"}
